# Update res_bus/vm_pu.xlsx values for the Case_3_145 (380 kV) scenario.
# Slack bus voltage setpoint (column B) changes from 1.05 to 1.02 p.u.,
# which in turn changes the computed per-unit voltages in columns C-F and I-N
# (column G stays at 1, column H has no data) for every data row (2-25).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.053654235886896
$ws.Range("D2").Value = 1.051906203679656
$ws.Range("E2").Value = 1.059927366696936
$ws.Range("F2").Value = 1.06990769793937
$ws.Range("I2").Value = 1.04029523161639
$ws.Range("J2").Value = 1.05867044515653
$ws.Range("K2").Value = 1.054656408621604
$ws.Range("L2").Value = 1.06265553713476
$ws.Range("M2").Value = 1.072608950954044
$ws.Range("N2").Value = 1.06017387861985

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.05527350981188
$ws.Range("D3").Value = 1.053134889736492
$ws.Range("E3").Value = 1.061378307529421
$ws.Range("F3").Value = 1.071463493711631
$ws.Range("I3").Value = 1.040675615136836
$ws.Range("J3").Value = 1.059937440632246
$ws.Range("K3").Value = 1.05569676588293
$ws.Range("L3").Value = 1.063919194974583
$ws.Range("M3").Value = 1.073979176030893
$ws.Range("N3").Value = 1.061442673374468

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.056319374806056
$ws.Range("D4").Value = 1.053928016714577
$ws.Range("E4").Value = 1.062315635471433
$ws.Range("F4").Value = 1.072468801277535
$ws.Range("I4").Value = 1.0409195339805
$ws.Range("J4").Value = 1.060755004938172
$ws.Range("K4").Value = 1.056367473026695
$ws.Range("L4").Value = 1.064734822928517
$ws.Range("M4").Value = 1.074863908668226
$ws.Range("N4").Value = 1.062261398715485

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.056758608553565
$ws.Range("D5").Value = 1.054260995605208
$ws.Range("E5").Value = 1.062709330564589
$ws.Range("F5").Value = 1.07289110700024
$ws.Range("I5").Value = 1.041021549938661
$ws.Range("J5").Value = 1.061098174148213
$ws.Range("K5").Value = 1.056648852485943
$ws.Range("L5").Value = 1.065077231327845
$ws.Range("M5").Value = 1.075235404680281
$ws.Range("N5").Value = 1.062605055265164

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.056832331862872
$ws.Range("D6").Value = 1.054316877934758
$ws.Range("E6").Value = 1.062775412971551
$ws.Range("F6").Value = 1.072961995149174
$ws.Range("I6").Value = 1.041038648010345
$ws.Range("J6").Value = 1.061155762631217
$ws.Range("K6").Value = 1.056696063166006
$ws.Range("L6").Value = 1.065134695182215
$ws.Range("M6").Value = 1.075297754596832
$ws.Range("N6").Value = 1.062662725530421

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.056325245614648
$ws.Range("D7").Value = 1.053932467760145
$ws.Range("E7").Value = 1.062320897438132
$ws.Range("F7").Value = 1.072474445414421
$ws.Range("I7").Value = 1.040920899191741
$ws.Range("J7").Value = 1.060759592472598
$ws.Range("K7").Value = 1.056371235125795
$ws.Range("L7").Value = 1.064739400088104
$ws.Range("M7").Value = 1.074868874355833
$ws.Range("N7").Value = 1.062265992764735

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.05420187945727
$ws.Range("D8").Value = 1.052321844506075
$ws.Range("E8").Value = 1.060418040227783
$ws.Range("F8").Value = 1.070433780782111
$ws.Range("I8").Value = 1.040424244206319
$ws.Range("J8").Value = 1.059099106198438
$ws.Range("K8").Value = 1.055008518236636
$ws.Range("L8").Value = 1.063083023826536
$ws.Range("M8").Value = 1.073072422048358
$ws.Range("N8").Value = 1.060603148409606

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.050445122874364
$ws.Range("D9").Value = 1.049468748108061
$ws.Range("E9").Value = 1.057052893264916
$ws.Range("F9").Value = 1.066826782980914
$ws.Range("I9").Value = 1.039531991275264
$ws.Range("J9").Value = 1.056155408523367
$ws.Range("K9").Value = 1.052588006141697
$ws.Range("L9").Value = 1.060148289355267
$ws.Range("M9").Value = 1.069891973706075
$ws.Range("N9").Value = 1.057655270346189

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.047929812552859
$ws.Range("D10").Value = 1.04755618738903
$ws.Range("E10").Value = 1.054800829532844
$ws.Range("F10").Value = 1.064414113794188
$ws.Range("I10").Value = 1.038925499071536
$ws.Range("J10").Value = 1.054180544969527
$ws.Range("K10").Value = 1.050960998886194
$ws.Range("L10").Value = 1.058180572955171
$ws.Range("M10").Value = 1.067761172073462
$ws.Range("N10").Value = 1.055677602259583

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.046837948392199
$ws.Range("D11").Value = 1.046725441377653
$ws.Range("E11").Value = 1.053823497643691
$ws.Range("F11").Value = 1.063367379731665
$ws.Range("I11").Value = 1.038660077864401
$ws.Range("J11").Value = 1.053322356253139
$ws.Range("K11").Value = 1.050253234087366
$ws.Range("L11").Value = 1.057325760231792
$ws.Range("M11").Value = 1.066835907189416
$ws.Range("N11").Value = 1.054818194816782

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.046431960564188
$ws.Range("D12").Value = 1.046416467101457
$ws.Range("E12").Value = 1.053460136457913
$ws.Range("F12").Value = 1.062978260051878
$ws.Range("I12").Value = 1.038561063642089
$ws.Range("J12").Value = 1.053003117468278
$ws.Range("K12").Value = 1.049989841078048
$ws.Range("L12").Value = 1.057007817979259
$ws.Range("M12").Value = 1.066491820096714
$ws.Range("N12").Value = 1.054498502676229

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.04651906557807
$ws.Range("D13").Value = 1.046482761288653
$ws.Range("E13").Value = 1.053538094060584
$ws.Range("F13").Value = 1.063061742018548
$ws.Range("I13").Value = 1.038582321814312
$ws.Range("J13").Value = 1.053071616662346
$ws.Range("K13").Value = 1.050046362419734
$ws.Range("L13").Value = 1.05707603713316
$ws.Range("M13").Value = 1.066565646319772
$ws.Range("N13").Value = 1.054567099147009

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.046804397969801
$ws.Range("D14").Value = 1.046699909640074
$ws.Range("E14").Value = 1.053793469019992
$ws.Range("F14").Value = 1.063335221486978
$ws.Range("I14").Value = 1.038651902006827
$ws.Range("J14").Value = 1.053295977504654
$ws.Range("K14").Value = 1.050231472136831
$ws.Range("L14").Value = 1.057299487785636
$ws.Range("M14").Value = 1.066807473111449
$ws.Range("N14").Value = 1.054791778607447

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.04698014452957
$ws.Range("D15").Value = 1.046833648885386
$ws.Range("E15").Value = 1.053950769054378
$ws.Range("F15").Value = 1.063503679012879
$ws.Range("I15").Value = 1.038694716245596
$ws.Range("J15").Value = 1.053434151151255
$ws.Range("K15").Value = 1.050345458162358
$ws.Range("L15").Value = 1.057437106279601
$ws.Range("M15").Value = 1.066956416945838
$ws.Range("N15").Value = 1.054930148476476

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.048002217648857
$ws.Range("D16").Value = 1.047611265884652
$ws.Range("E16").Value = 1.054865645182185
$ws.Range("F16").Value = 1.064483538355468
$ws.Range("I16").Value = 1.03894305479375
$ws.Range("J16").Value = 1.05423743484602
$ws.Range("K16").Value = 1.051007901550958
$ws.Range("L16").Value = 1.058237244727696
$ws.Range("M16").Value = 1.067822522951031
$ws.Range("N16").Value = 1.055734572926225

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.048642600596481
$ws.Range("D17").Value = 1.04809834392441
$ws.Range("E17").Value = 1.055438933834502
$ws.Range("F17").Value = 1.065097626582608
$ws.Range("I17").Value = 1.039098077356344
$ws.Range("J17").Value = 1.054740487556759
$ws.Range("K17").Value = 1.051422556471589
$ws.Range("L17").Value = 1.058738400002621
$ws.Range("M17").Value = 1.06836510168382
$ws.Range("N17").Value = 1.056238340029515

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.049015863718559
$ws.Range("D18").Value = 1.048382198535827
$ws.Range("E18").Value = 1.055773114387029
$ws.Range("F18").Value = 1.065455618488606
$ws.Range("I18").Value = 1.03918822887897
$ws.Range("J18").Value = 1.055033615121642
$ws.Range("K18").Value = 1.051664103673138
$ws.Range("L18").Value = 1.059030447940222
$ws.Range("M18").Value = 1.068681327162692
$ws.Range("N18").Value = 1.056531883869162

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.049143092832266
$ws.Range("D19").Value = 1.048478943510754
$ws.Range("E19").Value = 1.055887026169315
$ws.Range("F19").Value = 1.065577651661572
$ws.Range("I19").Value = 1.039218922436181
$ws.Range("J19").Value = 1.055133514370525
$ws.Range("K19").Value = 1.051746412004348
$ws.Range("L19").Value = 1.05912998365463
$ws.Range("M19").Value = 1.068789109513481
$ws.Range("N19").Value = 1.056631924986437

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.048573920696614
$ws.Range("D20").Value = 1.048046110955993
$ws.Range("E20").Value = 1.055377447080377
$ws.Range("F20").Value = 1.065031761012706
$ws.Range("I20").Value = 1.039081472907681
$ws.Range("J20").Value = 1.054686545299915
$ws.Range("K20").Value = 1.051378100458221
$ws.Range("L20").Value = 1.058684658515989
$ws.Range("M20").Value = 1.068306914207391
$ws.Range("N20").Value = 1.056184321168478

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.04672038641088
$ws.Range("D21").Value = 1.046635975923583
$ws.Range("E21").Value = 1.05371827683284
$ws.Range("F21").Value = 1.063254697410578
$ws.Range("I21").Value = 1.038631424132933
$ws.Range("J21").Value = 1.053229921860652
$ws.Range("K21").Value = 1.050176975740142
$ws.Range("L21").Value = 1.057233698993272
$ws.Range("M21").Value = 1.066736272339564
$ws.Range("N21").Value = 1.054725629156855

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.045552554586788
$ws.Range("D22").Value = 1.045747059866096
$ws.Range("E22").Value = 1.052673138088313
$ws.Range("F22").Value = 1.062135553093548
$ws.Range("I22").Value = 1.03834600016976
$ws.Range("J22").Value = 1.052311364939847
$ws.Range("K22").Value = 1.049418897749278
$ws.Range("L22").Value = 1.056318949458878
$ws.Range("M22").Value = 1.065746413679867
$ws.Range("N22").Value = 1.053805767779863

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.046171879932277
$ws.Range("D23").Value = 1.046218512705331
$ws.Range("E23").Value = 1.053227374377741
$ws.Range("F23").Value = 1.062729010100439
$ws.Range("I23").Value = 1.038497543172958
$ws.Range("J23").Value = 1.052798570387321
$ws.Range("K23").Value = 1.049821045181501
$ws.Range("L23").Value = 1.056804113099869
$ws.Range("M23").Value = 1.066271381006102
$ws.Range("N23").Value = 1.054293665114954

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.048604954994896
$ws.Range("D24").Value = 1.048069713557996
$ws.Range("E24").Value = 1.055405230943888
$ws.Range("F24").Value = 1.065061523432406
$ws.Range("I24").Value = 1.039088976579234
$ws.Range("J24").Value = 1.054710920394164
$ws.Range("K24").Value = 1.051398189187415
$ws.Range("L24").Value = 1.058708942807049
$ws.Range("M24").Value = 1.068333207400628
$ws.Range("N24").Value = 1.056208730878158

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.051418191578537
$ws.Range("D25").Value = 1.050208160233248
$ws.Range("E25").Value = 1.0579243477225
$ws.Range("F25").Value = 1.067760648495534
$ws.Range("I25").Value = 1.039764702016885
$ws.Range("J25").Value = 1.056918573902545
$ws.Range("K25").Value = 1.053216086730567
$ws.Range("L25").Value = 1.060908932100514
$ws.Range("M25").Value = 1.070716009723065
$ws.Range("N25").Value = 1.058419519507741
